$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.294.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.42%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.516.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.10%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.20%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.21%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.584"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.46%  "

# Row 8
$ws.Range("E8").Value = "  +0.09%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.70%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.72%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0808"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.23%  "

# Row 12
$ws.Range("E12").Value = "  -1.55%  "

# Row 13
$ws.Range("E13").Value = "  -2.63%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.901.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.10%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.523.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.60%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.37%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.855"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.32%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.340.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.36%  "

# Row 19
$ws.Range("E19").Value = "  -0.26%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0970"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.80%  "

# Row 21
$ws.Range("E21").Value = "  -2.96%  "

# Row 22
$ws.Range("E22").Value = "  -1.55%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.63%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.35%  "

# Row 25
$ws.Range("E25").Value = "  -5.30%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.71%  "

# Row 27
$ws.Range("E27").Value = "  +0.42%  "

# Row 28
$ws.Range("E28").Value = "  +9.60%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.06%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.08%  "

# Row 31
$ws.Range("E31").Value = "  -1.87%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "154.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.43%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.28%  "

# Row 34
$ws.Range("E34").Value = "  -2.57%  "

# Row 35
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.62"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.23%  "

# Row 36
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.99%  "

# Row 37
$ws.Range("E37").Value = "  -5.14%  "

# Row 38
$ws.Range("E38").Value = "  +1.15%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.88%  "

# Row 41
$ws.Range("E41").Value = "  -1.51%  "

# Row 42
$ws.Range("E42").Value = "  -1.40%  "

# Row 44
$ws.Range("E44").Value = "  -2.76%  "

# Row 45
$ws.Range("E45").Value = "  -4.00%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.031.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.14%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.89%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.76%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.762.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.13%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "101.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.98%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.187"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.14%  "
